# ---------------------------------------------------------------------------
# C5-PowerPoint.pptx edit
#
# Per the authoritative diff this commit makes two logical changes:
#
#   1. The table on slide 6 (the "SOURCES OF FINANCE" table) gets a new
#      table style applied: {4FACAF59-0A37-4440-80EB-85A256DCEC40} ->
#      {95C0DE45-82D8-403E-BCA8-E0CB752FFECF}.
#
#   2. The deck's theme palette is switched from the "Integral" colour
#      scheme to the standard "Office Theme" colour scheme (the colours
#      that, before the edit, were only used by the Notes Master's theme
#      part). We reproduce that by repointing every theme colour slot
#      (dark1/light1/dark2/light2/accent1-6/hyperlink/followed-hyperlink)
#      to the "Office Theme" RGB values.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 ---------------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{95C0DE45-82D8-403E-BCA8-E0CB752FFECF}")

# --- 2. Theme colour scheme -> "Office Theme" palette -----------------------
# RGB() packs colours as 0x00BBGGRR, so compute the decimal values from the
# target hex colours (dk1=000000, lt1=FFFFFF, dk2=44546A, lt2=E7E6E6,
# accent1=5B9BD5, accent2=ED7D31, accent3=A5A5A5, accent4=FFC000,
# accent5=4472C4, accent6=70AD47, hlink=0563C1, folHlink=954F72).
$themeColors = $p.Slides.Item(1).ThemeColorScheme
$themeColors.Item(1).RGB  = 0          # Dark 1      -> #000000
$themeColors.Item(2).RGB  = 16777215   # Light 1     -> #FFFFFF
$themeColors.Item(3).RGB  = 6968388    # Dark 2      -> #44546A
$themeColors.Item(4).RGB  = 15132391   # Light 2     -> #E7E6E6
$themeColors.Item(5).RGB  = 13998939   # Accent 1    -> #5B9BD5
$themeColors.Item(6).RGB  = 3243501    # Accent 2    -> #ED7D31
$themeColors.Item(7).RGB  = 10855845   # Accent 3    -> #A5A5A5
$themeColors.Item(8).RGB  = 49407      # Accent 4    -> #FFC000
$themeColors.Item(9).RGB  = 12874308   # Accent 5    -> #4472C4
$themeColors.Item(10).RGB = 4697456    # Accent 6    -> #70AD47
$themeColors.Item(11).RGB = 12673797   # Hyperlink   -> #0563C1
$themeColors.Item(12).RGB = 7491477    # Followed Hl -> #954F72
